$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '244.23'
$ws.Range('D2').Style = 'Normal'

$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '-0.92%'
$ws.Range('E2').Style = 'Normal'

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '27.20'
$ws.Range('D3').Style = 'Normal'

$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '4.05%'
$ws.Range('E3').Style = 'Normal'

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.157'
$ws.Range('D4').Style = 'Normal'

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '1.17%'
$ws.Range('E4').Style = 'Normal'

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.05635'
$ws.Range('D5').Style = 'Normal'

$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '0.67%'
$ws.Range('E5').Style = 'Normal'

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '6.478'
$ws.Range('D6').Style = 'Normal'

$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '-0.38%'
$ws.Range('E6').Style = 'Normal'

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.8161'
$ws.Range('D7').Style = 'Normal'

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '0.64%'
$ws.Range('E7').Style = 'Normal'

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.8318'
$ws.Range('D8').Style = 'Normal'

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '-1.20%'
$ws.Range('E8').Style = 'Normal'

$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '-1.33%'
$ws.Range('E9').Style = 'Normal'

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06922'
$ws.Range('D10').Style = 'Normal'

$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '-1.21%'
$ws.Range('E10').Style = 'Normal'

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.02926'
$ws.Range('D11').Style = 'Normal'

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '2.73%'
$ws.Range('E11').Style = 'Normal'

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.09398'
$ws.Range('D12').Style = 'Normal'

$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '-0.07%'
$ws.Range('E12').Style = 'Normal'

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.001514'
$ws.Range('D13').Style = 'Normal'

$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '0.10%'
$ws.Range('E13').Style = 'Normal'

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.04236'
$ws.Range('D14').Style = 'Normal'

$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '-9.61%'
$ws.Range('E14').Style = 'Normal'

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0005948'
$ws.Range('D15').Style = 'Normal'

$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '-1.03%'
$ws.Range('E15').Style = 'Normal'

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.006148'
$ws.Range('D16').Style = 'Normal'

$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '0.46%'
$ws.Range('E16').Style = 'Normal'

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.533'
$ws.Range('D17').Style = 'Normal'

$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '-0.64%'
$ws.Range('E17').Style = 'Normal'

$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '-0.53%'
$ws.Range('E18').Style = 'Normal'

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.309'
$ws.Range('D19').Style = 'Normal'

$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '9.03%'
$ws.Range('E19').Style = 'Normal'

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.03111'
$ws.Range('D21').Style = 'Normal'

$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '-3.65%'
$ws.Range('E21').Style = 'Normal'

$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '-2.15%'
$ws.Range('E22').Style = 'Normal'

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.738'
$ws.Range('D23').Style = 'Normal'

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '-0.11%'
$ws.Range('E23').Style = 'Normal'

$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '-0.09%'
$ws.Range('E24').Style = 'Normal'

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.001224'
$ws.Range('D25').Style = 'Normal'

$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '-1.74%'
$ws.Range('E25').Style = 'Normal'

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.004471'
$ws.Range('D26').Style = 'Normal'

$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '-2.93%'
$ws.Range('E26').Style = 'Normal'

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.00009797'
$ws.Range('D27').Style = 'Normal'

$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '2.06%'
$ws.Range('E27').Style = 'Normal'

$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '-0.49%'
$ws.Range('E28').Style = 'Normal'

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.03648'
$ws.Range('D40').Style = 'Normal'

$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '-0.18%'
$ws.Range('E40').Style = 'Normal'

$ws.Range('B41').Value = 'KickToken'

$ws.Range('C41').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.006068'
$ws.Range('D41').Style = 'Normal'

$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '-0.68%'
$ws.Range('E41').Style = 'Normal'

$ws.Range('B42').Value = 'BKEXToken'

$ws.Range('C42').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1052'
$ws.Range('D42').Style = 'Normal'

$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '-0.08%'
$ws.Range('E42').Style = 'Normal'

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.002299'
$ws.Range('D43').Style = 'Normal'

$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '-8.02%'
$ws.Range('E43').Style = 'Normal'

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.008221'
$ws.Range('D44').Style = 'Normal'

$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '-5.63%'
$ws.Range('E44').Style = 'Normal'

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00005388'
$ws.Range('D45').Style = 'Normal'

$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '1.84%'
$ws.Range('E45').Style = 'Normal'

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.1299'
$ws.Range('D47').Style = 'Normal'

$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '-2.14%'
$ws.Range('E47').Style = 'Normal'

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.004725'
$ws.Range('D48').Style = 'Normal'

$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '130.20%'
$ws.Range('E48').Style = 'Normal'
